$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Road Usage data (A: road id, B: usage count) - final version before presentation
$data = @(
    @(102, 2855),
    @(98, 3311),
    @(151, 3311),
    @(66, 782),
    @(162, 770),
    @(161, 758),
    @(63, 1484),
    @(62, 1484),
    @(79, 4842),
    @(64, 691),
    @(51, 854),
    @(53, 750),
    @(52, 883),
    @(55, 883),
    @(54, 883),
    @(57, 633),
    @(58, 633),
    @(61, 633),
    @(60, 633),
    @(116, 2128),
    @(114, 2128),
    @(115, 2020),
    @(159, 323),
    @(117, 2203),
    @(122, 3049),
    @(152, 3028),
    @(123, 2369),
    @(118, 2311),
    @(121, 2083),
    @(78, 3418),
    @(77, 2238),
    @(91, 2238),
    @(90, 2376),
    @(93, 2376),
    @(94, 2305),
    @(95, 2305),
    @(96, 2331),
    @(92, 2423),
    @(97, 2423),
    @(99, 2423),
    @(100, 2433),
    @(101, 2435),
    @(86, 2964),
    @(87, 2964),
    @(124, 390),
    @(125, 173),
    @(126, 173),
    @(68, 173),
    @(72, 173),
    @(70, 173),
    @(71, 179),
    @(69, 179),
    @(80, 4627),
    @(149, 4425),
    @(148, 4425),
    @(49, 4389),
    @(46, 998),
    @(85, 3030),
    @(84, 3030),
    @(81, 3191),
    @(82, 3840),
    @(83, 3840),
    @(48, 2830),
    @(76, 346),
    @(73, 318),
    @(75, 318),
    @(47, 1799),
    @(133, 214),
    @(128, 404),
    @(74, 239),
    @(38, 5374),
    @(39, 5060),
    @(42, 5060),
    @(43, 4928),
    @(40, 4928),
    @(41, 4787),
    @(44, 4311),
    @(45, 4287),
    @(14, 4287),
    @(13, 4316),
    @(12, 4316),
    @(135, 211),
    @(28, 774),
    @(27, 769),
    @(150, 168),
    @(142, 168),
    @(127, 302),
    @(120, 477),
    @(113, 990),
    @(105, 990),
    @(106, 856),
    @(129, 167),
    @(103, 568),
    @(158, 36),
    @(11, 3100),
    @(16, 107),
    @(107, 274),
    @(139, 3097),
    @(23, 3096),
    @(108, 39),
    @(109, 39),
    @(110, 39),
    @(111, 39),
    @(25, 11),
    @(24, 9),
    @(37, 42),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
